# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row -> new F value (applies identically to both sheets)
$updates = @{
    4  = 871
    5  = 39
    6  = 339
    7  = 10592
    8  = 178
    13 = 137
    16 = 38
    19 = 298
    20 = 996
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
